# Update NATMI LR-pair (Vcan-Cd44) TPM-derived statistics with newly
# recomputed values (ligand/receptor expression, specificity, and edge
# weight columns E-J, M-T) for every Sending x Target cluster combination.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 4.204118999999999
$ws.Cells.Item(2, 8).Value = 12.612357
$ws.Cells.Item(2, 9).Value = 0.01983154129720676
$ws.Cells.Item(2, 10).Value = 0.01983154129720676
$ws.Cells.Item(2, 13).Value = 24.576554
$ws.Cells.Item(2, 14).Value = 73.729662
$ws.Cells.Item(2, 15).Value = 0.07553767049546639
$ws.Cells.Item(2, 16).Value = 0.07553767049546638
$ws.Cells.Item(2, 17).Value = 103.322757625926
$ws.Cells.Item(2, 18).Value = 929.904818633334
$ws.Cells.Item(2, 19).Value = 0.001498028431925638
$ws.Cells.Item(2, 20).Value = 0.001498028431925638

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 4.204118999999999
$ws.Cells.Item(3, 8).Value = 12.612357
$ws.Cells.Item(3, 9).Value = 0.01983154129720676
$ws.Cells.Item(3, 10).Value = 0.01983154129720676
$ws.Cells.Item(3, 15).Value = 0.359764849016532
$ws.Cells.Item(3, 16).Value = 0.359764849016532
$ws.Cells.Item(3, 17).Value = 492.0974667797569
$ws.Cells.Item(3, 18).Value = 4428.877201017813
$ws.Cells.Item(3, 19).Value = 0.007134691460554708
$ws.Cells.Item(3, 20).Value = 0.007134691460554709

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 4.204118999999999
$ws.Cells.Item(4, 8).Value = 12.612357
$ws.Cells.Item(4, 9).Value = 0.01983154129720676
$ws.Cells.Item(4, 10).Value = 0.01983154129720676
$ws.Cells.Item(4, 13).Value = 55.68784966666667
$ws.Cells.Item(4, 14).Value = 167.063549
$ws.Cells.Item(4, 15).Value = 0.1711603033819035
$ws.Cells.Item(4, 16).Value = 0.1711603033819035
$ws.Cells.Item(4, 17).Value = 234.118346852777
$ws.Cells.Item(4, 18).Value = 2107.065121674993
$ws.Cells.Item(4, 19).Value = 0.003394372624960657
$ws.Cells.Item(4, 20).Value = 0.003394372624960658

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4.204118999999999
$ws.Cells.Item(5, 8).Value = 12.612357
$ws.Cells.Item(5, 9).Value = 0.01983154129720676
$ws.Cells.Item(5, 10).Value = 0.01983154129720676
$ws.Cells.Item(5, 13).Value = 128.0392633333333
$ws.Cells.Item(5, 14).Value = 384.11779
$ws.Cells.Item(5, 15).Value = 0.3935371771060981
$ws.Cells.Item(5, 16).Value = 0.3935371771060981
$ws.Cells.Item(5, 17).Value = 538.29229972567
$ws.Cells.Item(5, 18).Value = 4844.63069753103
$ws.Cells.Item(5, 19).Value = 0.007804448779765754
$ws.Cells.Item(5, 20).Value = 0.007804448779765755

$ws.Cells.Item(6, 9).Value = 0.8539093107807857
$ws.Cells.Item(6, 10).Value = 0.8539093107807858
$ws.Cells.Item(6, 13).Value = 24.576554
$ws.Cells.Item(6, 14).Value = 73.729662
$ws.Cells.Item(6, 15).Value = 0.07553767049546639
$ws.Cells.Item(6, 16).Value = 0.07553767049546638
$ws.Cells.Item(6, 17).Value = 4448.88591512307
$ws.Cells.Item(6, 18).Value = 40039.97323610763
$ws.Cells.Item(6, 19).Value = 0.0645023201507698
$ws.Cells.Item(6, 20).Value = 0.06450232015076979

$ws.Cells.Item(7, 9).Value = 0.8539093107807857
$ws.Cells.Item(7, 10).Value = 0.8539093107807858
$ws.Cells.Item(7, 15).Value = 0.359764849016532
$ws.Cells.Item(7, 16).Value = 0.359764849016532
$ws.Cells.Item(7, 19).Value = 0.3072065542668603
$ws.Cells.Item(7, 20).Value = 0.3072065542668603

$ws.Cells.Item(8, 9).Value = 0.8539093107807857
$ws.Cells.Item(8, 10).Value = 0.8539093107807858
$ws.Cells.Item(8, 13).Value = 55.68784966666667
$ws.Cells.Item(8, 14).Value = 167.063549
$ws.Cells.Item(8, 15).Value = 0.1711603033819035
$ws.Cells.Item(8, 16).Value = 0.1711603033819035
$ws.Cells.Item(8, 17).Value = 10080.70089995222
$ws.Cells.Item(8, 18).Value = 90726.30809956997
$ws.Cells.Item(8, 19).Value = 0.1461553766938714
$ws.Cells.Item(8, 20).Value = 0.1461553766938714

$ws.Cells.Item(9, 9).Value = 0.8539093107807857
$ws.Cells.Item(9, 10).Value = 0.8539093107807858
$ws.Cells.Item(9, 13).Value = 128.0392633333333
$ws.Cells.Item(9, 14).Value = 384.11779
$ws.Cells.Item(9, 15).Value = 0.3935371771060981
$ws.Cells.Item(9, 16).Value = 0.3935371771060981
$ws.Cells.Item(9, 17).Value = 23177.86599481768
$ws.Cells.Item(9, 18).Value = 208600.7939533591
$ws.Cells.Item(9, 19).Value = 0.3360450596692842
$ws.Cells.Item(9, 20).Value = 0.3360450596692843

$ws.Cells.Item(10, 7).Value = 26.057747
$ws.Cells.Item(10, 8).Value = 78.173241
$ws.Cells.Item(10, 9).Value = 0.1229188055196976
$ws.Cells.Item(10, 10).Value = 0.1229188055196976
$ws.Cells.Item(10, 13).Value = 24.576554
$ws.Cells.Item(10, 14).Value = 73.729662
$ws.Cells.Item(10, 15).Value = 0.07553767049546639
$ws.Cells.Item(10, 16).Value = 0.07553767049546638
$ws.Cells.Item(10, 17).Value = 640.409626263838
$ws.Cells.Item(10, 18).Value = 5763.686636374543
$ws.Cells.Item(10, 19).Value = 0.009285000229043234
$ws.Cells.Item(10, 20).Value = 0.009285000229043232

$ws.Cells.Item(11, 7).Value = 26.057747
$ws.Cells.Item(11, 8).Value = 78.173241
$ws.Cells.Item(11, 9).Value = 0.1229188055196976
$ws.Cells.Item(11, 10).Value = 0.1229188055196976
$ws.Cells.Item(11, 15).Value = 0.359764849016532
$ws.Cells.Item(11, 16).Value = 0.359764849016532
$ws.Cells.Item(11, 17).Value = 3050.092371002775
$ws.Cells.Item(11, 18).Value = 27450.83133902497
$ws.Cells.Item(11, 19).Value = 0.04422186550908647
$ws.Cells.Item(11, 20).Value = 0.04422186550908647

$ws.Cells.Item(12, 7).Value = 26.057747
$ws.Cells.Item(12, 8).Value = 78.173241
$ws.Cells.Item(12, 9).Value = 0.1229188055196976
$ws.Cells.Item(12, 10).Value = 0.1229188055196976
$ws.Cells.Item(12, 13).Value = 55.68784966666667
$ws.Cells.Item(12, 14).Value = 167.063549
$ws.Cells.Item(12, 15).Value = 0.1711603033819035
$ws.Cells.Item(12, 16).Value = 0.1711603033819035
$ws.Cells.Item(12, 17).Value = 1451.099897588035
$ws.Cells.Item(12, 18).Value = 13059.89907829231
$ws.Cells.Item(12, 19).Value = 0.02103882004409264
$ws.Cells.Item(12, 20).Value = 0.02103882004409264

$ws.Cells.Item(13, 7).Value = 26.057747
$ws.Cells.Item(13, 8).Value = 78.173241
$ws.Cells.Item(13, 9).Value = 0.1229188055196976
$ws.Cells.Item(13, 10).Value = 0.1229188055196976
$ws.Cells.Item(13, 13).Value = 128.0392633333333
$ws.Cells.Item(13, 14).Value = 384.11779
$ws.Cells.Item(13, 15).Value = 0.3935371771060981
$ws.Cells.Item(13, 16).Value = 0.3935371771060981
$ws.Cells.Item(13, 17).Value = 3336.414730006377
$ws.Cells.Item(13, 18).Value = 30027.73257005739
$ws.Cells.Item(13, 19).Value = 0.04837311973747527
$ws.Cells.Item(13, 20).Value = 0.04837311973747527

$ws.Cells.Item(14, 7).Value = 0.7081243333333332
$ws.Cells.Item(14, 8).Value = 2.124373
$ws.Cells.Item(14, 9).Value = 0.003340342402309973
$ws.Cells.Item(14, 10).Value = 0.003340342402309974
$ws.Cells.Item(14, 13).Value = 24.576554
$ws.Cells.Item(14, 14).Value = 73.729662
$ws.Cells.Item(14, 15).Value = 0.07553767049546639
$ws.Cells.Item(14, 16).Value = 0.07553767049546638
$ws.Cells.Item(14, 17).Value = 17.40325591688067
$ws.Cells.Item(14, 18).Value = 156.629303251926
$ws.Cells.Item(14, 19).Value = 0.0002523216837277254
$ws.Cells.Item(14, 20).Value = 0.0002523216837277254

$ws.Cells.Item(15, 7).Value = 0.7081243333333332
$ws.Cells.Item(15, 8).Value = 2.124373
$ws.Cells.Item(15, 9).Value = 0.003340342402309973
$ws.Cells.Item(15, 10).Value = 0.003340342402309974
$ws.Cells.Item(15, 15).Value = 0.359764849016532
$ws.Cells.Item(15, 16).Value = 0.359764849016532
$ws.Cells.Item(15, 17).Value = 82.88685229852855
$ws.Cells.Item(15, 18).Value = 745.981670686757
$ws.Cells.Item(15, 19).Value = 0.001201737780030567
$ws.Cells.Item(15, 20).Value = 0.001201737780030568

$ws.Cells.Item(16, 7).Value = 0.7081243333333332
$ws.Cells.Item(16, 8).Value = 2.124373
$ws.Cells.Item(16, 9).Value = 0.003340342402309973
$ws.Cells.Item(16, 10).Value = 0.003340342402309974
$ws.Cells.Item(16, 13).Value = 55.68784966666667
$ws.Cells.Item(16, 14).Value = 167.063549
$ws.Cells.Item(16, 15).Value = 0.1711603033819035
$ws.Cells.Item(16, 16).Value = 0.1711603033819035
$ws.Cells.Item(16, 17).Value = 39.43392141997522
$ws.Cells.Item(16, 18).Value = 354.905292779777
$ws.Cells.Item(16, 19).Value = 0.0005717340189788114
$ws.Cells.Item(16, 20).Value = 0.0005717340189788115

$ws.Cells.Item(17, 7).Value = 0.7081243333333332
$ws.Cells.Item(17, 8).Value = 2.124373
$ws.Cells.Item(17, 9).Value = 0.003340342402309973
$ws.Cells.Item(17, 10).Value = 0.003340342402309974
$ws.Cells.Item(17, 13).Value = 128.0392633333333
$ws.Cells.Item(17, 14).Value = 384.11779
$ws.Cells.Item(17, 15).Value = 0.3935371771060981
$ws.Cells.Item(17, 16).Value = 0.3935371771060981
$ws.Cells.Item(17, 17).Value = 90.66771798840777
$ws.Cells.Item(17, 18).Value = 816.0094618956699
$ws.Cells.Item(17, 19).Value = 0.001314548919572869
$ws.Cells.Item(17, 20).Value = 0.001314548919572869

